$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 176.7
$ws.Range("I5").Value = 179.4
$ws.Range("J5").Value = 174
$ws.Range("K5").Value = 179.4
$ws.Range("L5").Value = 174
$ws.Range("M5").Value = -64.40000000000001
$ws.Range("N5").Value = -404
$ws.Range("H40").Value = 2096.1
$ws.Range("I40").Value = 1992.2
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 1992.2
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -1817.2
$ws.Range("N40").Value = -2550
$ws.Range("H69").Value = 3301.88
$ws.Range("I69").Value = 3072
$ws.Range("J69").Value = 3374.4736
$ws.Range("K69").Value = 9216
$ws.Range("L69").Value = 10123.4208
$ws.Range("M69").Value = -8342
$ws.Range("N69").Value = -11871.4208
$ws.Range("H72").Value = 3301.88
$ws.Range("I72").Value = 3072
$ws.Range("J72").Value = 3374.4736
$ws.Range("K72").Value = 27648
$ws.Range("L72").Value = 30370.2624
$ws.Range("M72").Value = -23280
$ws.Range("N72").Value = -39106.2624
$ws.Range("H116").Value = 3094.818
$ws.Range("I116").Value = 2746.923
$ws.Range("K116").Value = 2746.923
$ws.Range("M116").Value = 695.0770000000002
$ws.Range("H135").Value = 1590.3182
$ws.Range("I135").Value = 610.3889
$ws.Range("K135").Value = 5493.5001
$ws.Range("M135").Value = -2958.5001
$ws.Range("H138").Value = 1807.5385
$ws.Range("I138").Value = 1319.1389
$ws.Range("J138").Value = 2906.4375
$ws.Range("K138").Value = 3957.4167
$ws.Range("L138").Value = 8719.3125
$ws.Range("M138").Value = 1182.5833
$ws.Range("N138").Value = -18999.3125

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1121.3
$ws.Range("I2").Value = 1183.3334
$ws.Range("J2").Value = 1028.25
$ws.Range("K2").Value = 1183.3334
$ws.Range("L2").Value = 1028.25
$ws.Range("M2").Value = -1070.3334
$ws.Range("N2").Value = -1254.25
$ws.Range("H19").Value = 17650
$ws.Range("I19").Value = 300
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -71
$ws.Range("H116").Value = 1121.3
$ws.Range("I116").Value = 1183.3334
$ws.Range("J116").Value = 1028.25
$ws.Range("K116").Value = 1183.3334
$ws.Range("L116").Value = 1028.25
$ws.Range("M116").Value = 1110.6666
$ws.Range("N116").Value = -5616.25

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1121.3
$ws.Range("I3").Value = 1183.3334
$ws.Range("J3").Value = 1028.25
$ws.Range("K3").Value = 1183.3334
$ws.Range("L3").Value = 1028.25
$ws.Range("M3").Value = -1069.3334
$ws.Range("N3").Value = -1256.25
$ws.Range("H22").Value = 3613.7
$ws.Range("I22").Value = 3777.5789
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 3777.5789
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -3604.5789
$ws.Range("N22").Value = -846
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 602.2
$ws.Range("I16").Value = 399.7143
$ws.Range("J16").Value = 1074.6666
$ws.Range("K16").Value = 399.7143
$ws.Range("L16").Value = 1074.6666
$ws.Range("M16").Value = -112.7143
$ws.Range("N16").Value = -1648.6666
$ws.Range("H31").Value = 4969.5
$ws.Range("I31").Value = 1174.9429
$ws.Range("J31").Value = 9253.678
$ws.Range("K31").Value = 1174.9429
$ws.Range("L31").Value = 9253.678
$ws.Range("M31").Value = -879.9429
$ws.Range("N31").Value = -9843.678
$ws.Range("H34").Value = 4969.5
$ws.Range("I34").Value = 1174.9429
$ws.Range("J34").Value = 9253.678
$ws.Range("K34").Value = 1174.9429
$ws.Range("L34").Value = 9253.678
$ws.Range("M34").Value = -972.9429
$ws.Range("N34").Value = -9657.678
$ws.Range("H105").Value = 1742.4166
$ws.Range("I105").Value = 1434.4445
$ws.Range("K105").Value = 1434.4445
$ws.Range("M105").Value = 312.5554999999999
$ws.Range("H113").Value = 602.2
$ws.Range("I113").Value = 399.7143
$ws.Range("J113").Value = 1074.6666
$ws.Range("K113").Value = 399.7143
$ws.Range("L113").Value = 1074.6666
$ws.Range("M113").Value = 1770.2857
$ws.Range("N113").Value = -5414.6666

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7974.643
$ws.Range("I122").Value = 531.125
$ws.Range("J122").Value = 17899.334
$ws.Range("K122").Value = 4780.125
$ws.Range("L122").Value = 161094.006
$ws.Range("M122").Value = -2330.125
$ws.Range("N122").Value = -165994.006
$ws.Range("H132").Value = 2616.2083
$ws.Range("I132").Value = 2608.3333
$ws.Range("J132").Value = 2618.8333
$ws.Range("K132").Value = 23474.9997
$ws.Range("L132").Value = 23569.4997
$ws.Range("M132").Value = -20944.9997
$ws.Range("N132").Value = -28629.4997

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 42451760
$ws.Range("I80").Value = 72715860
$ws.Range("J80").Value = 82017.60000000001
$ws.Range("K80").Value = 72715860
$ws.Range("L80").Value = 82017.60000000001
$ws.Range("M80").Value = -72714862
$ws.Range("N80").Value = -84013.60000000001
$ws.Range("H83").Value = 42451760
$ws.Range("I83").Value = 72715860
$ws.Range("J83").Value = 82017.60000000001
$ws.Range("K83").Value = 363579300
$ws.Range("L83").Value = 410088
$ws.Range("M83").Value = -363574308
$ws.Range("N83").Value = -420072
$ws.Range("H102").Value = 1713.8572
$ws.Range("I102").Value = 1666.1666
$ws.Range("K102").Value = 1666.1666
$ws.Range("M102").Value = -44.16660000000002
$ws.Range("H113").Value = 72608
$ws.Range("I113").Value = 111766.555
$ws.Range("J113").Value = 2122.6
$ws.Range("K113").Value = 111766.555
$ws.Range("L113").Value = 2122.6
$ws.Range("M113").Value = -109596.555
$ws.Range("N113").Value = -6462.6
$ws.Range("H126").Value = 3635.3333
$ws.Range("I126").Value = 4003
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 12009
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -9539
$ws.Range("N126").Value = -13640
$ws.Range("H136").Value = 20603.572
$ws.Range("J136").Value = 20603.572
$ws.Range("L136").Value = 61810.716
$ws.Range("N136").Value = -66910.716

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 659.2
$ws.Range("I46").Value = 599
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 599
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -411
$ws.Range("N46").Value = -1276
$ws.Range("H68").Value = 1027
$ws.Range("I68").Value = 685
$ws.Range("J68").Value = 1300.6
$ws.Range("K68").Value = 685
$ws.Range("L68").Value = 1300.6
$ws.Range("M68").Value = 64
$ws.Range("N68").Value = -2798.6
$ws.Range("H71").Value = 1027
$ws.Range("I71").Value = 685
$ws.Range("J71").Value = 1300.6
$ws.Range("K71").Value = 3425
$ws.Range("L71").Value = 6503
$ws.Range("M71").Value = 319
$ws.Range("N71").Value = -13991

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6046.6
$ws.Range("I81").Value = 8567.5
$ws.Range("J81").Value = 4366
$ws.Range("K81").Value = 17135
$ws.Range("L81").Value = 8732
$ws.Range("M81").Value = -16074
$ws.Range("N81").Value = -10854
$ws.Range("H84").Value = 6046.6
$ws.Range("I84").Value = 8567.5
$ws.Range("J84").Value = 4366
$ws.Range("K84").Value = 85675
$ws.Range("L84").Value = 43660
$ws.Range("M84").Value = -80371
$ws.Range("N84").Value = -54268
$ws.Range("H96").Value = 4897.5
$ws.Range("I96").Value = 3616
$ws.Range("J96").Value = 7033.3335
$ws.Range("K96").Value = 3616
$ws.Range("L96").Value = 7033.3335
$ws.Range("M96").Value = -2243
$ws.Range("N96").Value = -9779.333500000001
$ws.Range("H126").Value = 3898.3333
$ws.Range("I126").Value = 896.6667
$ws.Range("K126").Value = 2690.0001
$ws.Range("M126").Value = -220.0001000000002
$ws.Range("H132").Value = 1660.3096
$ws.Range("I132").Value = 1450
$ws.Range("J132").Value = 1914.8948
$ws.Range("K132").Value = 4350
$ws.Range("L132").Value = 5744.6844
$ws.Range("M132").Value = -1820
$ws.Range("N132").Value = -10804.6844
$ws.Range("H137").Value = 41333.332
$ws.Range("J137").Value = 41333.332
$ws.Range("L137").Value = 41333.332
$ws.Range("N137").Value = -51533.332
